$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (price) values that look numeric from Excel auto-converting
# strings like "553.23" into numbers -- force text entry, then restore the
# original (default/"General") cell style so the resulting style matches the source.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '62.741.26'
$ws.Range("E2").Value = '  -1.64%  '
$ws.Range("D3").Value = '2.675.20'
$ws.Range("E3").Value = '  -2.20%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '553.23'
$ws.Range("E5").Value = '  -1.96%  '
$ws.Range("D6").Value = '156.83'
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("E9").Value = '  -3.69%  '
$ws.Range("E10").Value = '  -3.01%  '
$ws.Range("D11").Value = '0.365'
$ws.Range("E11").Value = '  -3.77%  '
$ws.Range("E12").Value = '  -4.28%  '
$ws.Range("D13").Value = '3.148.24'
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").Value = '26.19'
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("D15").Value = '62.655.98'
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("E16").Value = '  -2.89%  '
$ws.Range("D17").Value = '2.678.44'
$ws.Range("E17").Value = '  -2.19%  '
$ws.Range("D18").Value = '11.72'
$ws.Range("E18").Value = '  -6.87%  '
$ws.Range("D19").Value = '4.58'
$ws.Range("E19").Value = '  -3.26%  '
$ws.Range("D20").Value = '343.38'
$ws.Range("E20").Value = '  -2.98%  '
$ws.Range("E21").Value = '  -6.25%  '
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").Value = '0.508'
$ws.Range("E23").Value = '  -2.57%  '
$ws.Range("D24").Value = '63.10'
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").Value = '8.10'
$ws.Range("E27").Value = '  -3.36%  '
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '1.38'
$ws.Range("E28").Value = '  +4.88%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0845'
$ws.Range("E29").Value = '  -6.81%  '
$ws.Range("D30").Value = '7.24'
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("E31").Value = '  -1.75%  '
$ws.Range("D32").Value = '163.58'
$ws.Range("E32").Value = '  -0.86%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '4.83'
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").Value = '1.45'
$ws.Range("E35").Value = '  -1.22%  '
$ws.Range("D36").Value = '19.36'
$ws.Range("E36").Value = '  -3.38%  '
$ws.Range("E37").Value = '  -1.52%  '
$ws.Range("D38").Value = '340.20'
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '6.12'
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").Value = '0.930'
$ws.Range("E40").Value = '  -4.67%  '
$ws.Range("E41").Value = '  -3.49%  '
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("D43").Value = '20.61'
$ws.Range("E43").Value = '  -6.24%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '20.08'
$ws.Range("E44").Value = '  -4.77%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '0.614'
$ws.Range("E45").Value = '  -1.90%  '
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").Value = '11.01'
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '0.0552'
$ws.Range("E48").Value = '  -5.47%  '
$ws.Range("D49").Value = '0.0968'
$ws.Range("E49").Value = '  -3.27%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '0.0240'
$ws.Range("E50").Value = '  -4.28%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '127.93'
$ws.Range("E51").Value = '  -3.39%  '

# Restore default styling on the price column (matches original unstyled cells)
$priceRange.Style = "Normal"

